# #3438 - insert field on excel
#
# Inserts a new "Journal Number" column into the Advance Payment Report
# header row (row 8), between the existing "Supplier Invoice Posting Date"
# column (K) and "Description" column (L, which becomes M).  Excel's normal
# EntireColumn.Insert() semantics push every column from L onward one slot
# to the right and carry the left-neighbour's formatting into the freshly
# inserted column, which matches the target column/row layout exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at L - everything from the old L (Description) onward
# shifts right by one, reusing column K's formatting for the new column.
$ws.Columns("L").Insert()

# Match the new column's width to its left neighbour (column K), since
# Insert() only copies cell-level formatting, not the column width.
$ws.Columns("L").ColumnWidth = $ws.Columns("K").ColumnWidth

# Populate the new header cell.
$ws.Range("L8").Value = "Journal Number"

# Reflect the author's resulting selection/scroll position on the sheet.
[void]$ws.Range("M11").Select()
